# Project Cleanup and added some missing translations
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# The existing ANSWER100QUESTIONS / LOSE1TIME rows (8 and 9) need to move down
# to rows 12 and 13 to make room for four new translation rows (MUSIC, SOUND,
# LANGUAGE, LEVELS) inserted right after the EXIT row (row 7).
$ws.Range("A8:D11").EntireRow.Insert()

# Column A (keys), top to bottom
$ws.Range("A8").Value = "MUSIC"
$ws.Range("A9").Value = "SOUND"
$ws.Range("A10").Value = "LANGUAGE"
$ws.Range("A11").Value = "LEVELS"

# Column B (Polish), top to bottom
$ws.Range("B8").Value = "Muzyka"
$ws.Range("B9").Value = "Dźwięk"
$ws.Range("B10").Value = "Język"
$ws.Range("B11").Value = "Poziomy"

# Column C (English), bottom to top
$ws.Range("C11").Value = "Levels"
$ws.Range("C10").Value = "Language"
$ws.Range("C9").Value = "Sound"
$ws.Range("C8").Value = "Music"

$ws.Range("C8").Select()
